$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Potential Arm Stuff")

# Insert a new column before the old "Supplier" column (C) to hold the
# supplier part number, shifting the rest of the table one column right.
$ws.Columns("C").Insert()
$ws.Range("C1").Value = "supplier number"

# Highlight the TinyDuino row (row 5) in yellow, like the rest of the
# spec rows that still need parts ordered.
$ws.Range("A5:I5").Interior.Color = 65535

# Replace the old (wrong) row 8 entry with the real motor / camera spec:
# an Adafruit Trinket mini MCU board, ordered from Digikey.
$ws.Range("A8").Value = 1500
$ws.Range("B8").Value = "adafruit"
$ws.Range("C8").Value = "1528-1020-ND"
$ws.Range("D8").Value = "digikey"
$ws.Range("E8").Value = "TRINKET MINI MCU BOARD 3.3V"
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 9.76
$ws.Range("H8").Formula = "=F8*G8"

$ws.Range("C8:E8").Font.Name = "Arial"
$ws.Range("C8:E8").Font.Size = 9
$ws.Range("C8:E8").Font.Color = 0

$ws.Rows("8").RowHeight = 15

# Remove the now-empty trailing row.
$ws.Rows("10").Delete()

$ws.Range("A3").Select()
